# Fruta / hortaliza, semanal
# Rotate the weekly price rows: row4 -> row2, row2 -> row3, row3 -> row4
# (Fecha, Volumen, Precio minimo/maximo/promedio, Origen, Precio $/Kg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return @{
        D = $ws.Range("D$row").Value2
        M = $ws.Range("M$row").Value2
        N = $ws.Range("N$row").Value2
        O = $ws.Range("O$row").Value2
        P = $ws.Range("P$row").Value2
        R = $ws.Range("R$row").Value2
        S = $ws.Range("S$row").Value2
    }
}

function Set-RowData($row, $data) {
    $ws.Range("D$row").Value = $data.D
    $ws.Range("M$row").Value = $data.M
    $ws.Range("N$row").Value = $data.N
    $ws.Range("O$row").Value = $data.O
    $ws.Range("P$row").Value = $data.P
    $ws.Range("R$row").Value = $data.R
    $ws.Range("S$row").Value = $data.S
}

$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row4 = Get-RowData 4

Set-RowData 2 $row4
Set-RowData 3 $row2
Set-RowData 4 $row3
